# Weekly update: insert a new daily price record for
# "Comercializadora del Agro de Limarí - Poroto granado" right after the most
# recent existing record (old row 72), pushing all subsequent rows (72-83)
# down by one (to 73-84), and growing the sheet dimension to A1:R84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; Excel shifts rows 72:83 down to 73:84
# and extends the used range automatically.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Range("A72").Value = 2
$ws.Range("B72").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44644
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = 100112030
$ws.Range("G72").Value = "Poroto granado"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 400
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 20000
$ws.Range("M72").Value = 19000
$ws.Range("N72").Value = "$/malla 25 kilos"
$ws.Range("O72").Value = "Provincia de Limarí"
$ws.Range("P72").Value = 760
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
